# Update gh-pages output (generated at 456a3b4)
#
# A new event ("上饶·ETI03动漫节") starting 2025-01-01 is added to the schedule.
# It slots in right where "南昌·萌卡动漫展" used to be listed; the 萌卡动漫展 row
# itself is pushed one row down (and its "想去人数" counter ticks from 1351 to
# 1354), which in turn pushes the following row ("九江·第二届异次元动漫嘉年华")
# down by one as well. Three other, earlier rows simply get their "想去人数"
# (column F) counters bumped up.
#
# This applies identically to two worksheets: "展览" (index 1) and
# "全部类型" (index 4) - the only difference being which row currently holds
# 萌卡动漫展 in each sheet.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function Set-TextCell($cell, [string]$text) {
    # Forces the value to be stored as plain text instead of letting Excel
    # auto-detect/convert date-looking strings (e.g. "2025-01-01") into date
    # serial numbers. ClearFormats() afterwards drops the temporary "@" text
    # number format again so the cell keeps its original (default) style.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Update-ExpoSheet($ws, [int]$mengkaRow) {
    # $mengkaRow = row number (1-based) that currently holds "南昌·萌卡动漫展"

    $cm04Row  = $mengkaRow - 4   # "南昌·CM04动漫游戏博览会"   F: 5180 -> 5201
    $suiyouRow = $mengkaRow - 3  # "南昌·岁酉山河..."          F: 113  -> 114
    $yunyaRow = $mengkaRow - 2   # "南昌·云芽动漫音乐嘉年华"   F: 5309 -> 5323

    # --- small numeric bumps on existing rows ---
    $ws.Cells.Item($cm04Row, 6).Value = 5201
    $ws.Cells.Item($suiyouRow, 6).Value = 114
    $ws.Cells.Item($yunyaRow, 6).Value = 5323

    # --- insert a new row right after the 萌卡动漫展 row, shifting rows below down ---
    $newRow = $mengkaRow + 1
    $ws.Rows($newRow).Insert()

    # The freshly inserted row's index cell (column A) needs the same
    # bold/centered/bordered style as the rest of that column; copy it from
    # a known-good row and then set its number explicitly.
    $ws.Range("A" + ($mengkaRow - 1)).Copy()
    $ws.Range("A" + $newRow).PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0
    $ws.Cells.Item($newRow, 1).Value = $mengkaRow

    # The row that got pushed down (previously "九江·..." etc.) keeps its own
    # content, but the sequential index in column A still needs to advance by
    # one to stay in step with its new row position.
    $ws.Cells.Item($newRow + 1, 1).Value = $newRow

    # --- move the (old) 萌卡动漫展 row content down into the freshly inserted row ---
    Set-TextCell $ws.Cells.Item($newRow, 2) "2025-01-01"
    Set-TextCell $ws.Cells.Item($newRow, 3) "南昌·萌卡动漫展"
    Set-TextCell $ws.Cells.Item($newRow, 4) "八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆"
    Set-TextCell $ws.Cells.Item($newRow, 5) "2025.01.01 09:00-01.03 17:00"
    $ws.Cells.Item($newRow, 6).Value = 1354
    $ws.Cells.Item($newRow, 7).Value = 65
    Set-TextCell $ws.Cells.Item($newRow, 8) "https://show.bilibili.com/platform/detail.html?id=93031"
    Set-TextCell $ws.Cells.Item($newRow, 9) "//i2.hdslb.com/bfs/openplatform/202409/HTlK8fN21727112669248.jpeg"

    # --- overwrite the original row with the new 上饶·ETI03动漫节 event ---
    # (column B, the start date, stays "2025-01-01" and is left untouched)
    Set-TextCell $ws.Cells.Item($mengkaRow, 3) "上饶·ETI03动漫节"
    Set-TextCell $ws.Cells.Item($mengkaRow, 4) "滨江东路与体育馆路交叉口西100米 力加体育综合运动中心"
    Set-TextCell $ws.Cells.Item($mengkaRow, 5) "2025.01.01 10:00-01.01 17:00"
    $ws.Cells.Item($mengkaRow, 6).Value = 0
    Set-TextCell $ws.Cells.Item($mengkaRow, 7) "不可售"
    Set-TextCell $ws.Cells.Item($mengkaRow, 8) "https://show.bilibili.com/platform/detail.html?id=93761"
    Set-TextCell $ws.Cells.Item($mengkaRow, 9) "//i0.hdslb.com/bfs/openplatform/202410/Ql2EXYVH1728884102415.jpeg"
}

# Sheet "展览" (index 1): 萌卡动漫展 currently sits on row 10
$wsExpo = $wb.Worksheets.Item(1)
Update-ExpoSheet $wsExpo 10

# Sheet "全部类型" (index 4): 萌卡动漫展 currently sits on row 11
$wsAll = $wb.Worksheets.Item(4)
Update-ExpoSheet $wsAll 11
